$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 910.34784
$ws.Range("J129").Value = 1042.3158
$ws.Range("L129").Value = 3126.9474
$ws.Range("N129").Value = -13126.9474

$ws.Range("H132").Value = 30498382
$ws.Range("I132").Value = 36754204
$ws.Range("J132").Value = 1250.625
$ws.Range("K132").Value = 110262612
$ws.Range("L132").Value = 3751.875
$ws.Range("M132").Value = -110260082
$ws.Range("N132").Value = -8811.875

$ws.Range("H135").Value = 2656.6309
$ws.Range("I135").Value = 2574.0527
$ws.Range("J135").Value = 3245
$ws.Range("K135").Value = 23166.4743
$ws.Range("L135").Value = 29205
$ws.Range("M135").Value = -20631.4743
$ws.Range("N135").Value = -34275

$ws.Range("H137").Value = 265939.75
$ws.Range("I137").Value = 332057.25
$ws.Range("J137").Value = 1469.7778
$ws.Range("K137").Value = 996171.75
$ws.Range("L137").Value = 4409.3334
$ws.Range("M137").Value = -993621.75
$ws.Range("N137").Value = -9509.3334

$ws.Range("H138").Value = 1656.6086
$ws.Range("I138").Value = 1049
$ws.Range("J138").Value = 2412.4146
$ws.Range("K138").Value = 3147
$ws.Range("L138").Value = 7237.2438
$ws.Range("M138").Value = 1993
$ws.Range("N138").Value = -17517.2438

$ws.Range("H140").Value = 34199
$ws.Range("J140").Value = 34199
$ws.Range("L140").Value = 34199
$ws.Range("N140").Value = -44559

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4659.6978
$ws.Range("I32").Value = 4378.141
$ws.Range("J32").Value = 5879.778
$ws.Range("K32").Value = 4378.141
$ws.Range("L32").Value = 5879.778
$ws.Range("M32").Value = -4091.141
$ws.Range("N32").Value = -6453.778

$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""  # cell removed entirely in target (was -150768)

$ws.Range("H132").Value = 2452008.2
$ws.Range("I132").Value = 2605225.2
$ws.Range("J132").Value = 537.6667
$ws.Range("K132").Value = 7815675.600000001
$ws.Range("L132").Value = 1613.0001
$ws.Range("M132").Value = -7813145.600000001
$ws.Range("N132").Value = -6673.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25680542
$ws.Range("I134").Value = 33384106
$ws.Range("J134").Value = 1996.6666
$ws.Range("K134").Value = 100152318
$ws.Range("L134").Value = 5989.9998
$ws.Range("M134").Value = -100149783
$ws.Range("N134").Value = -11059.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 23264396
$ws.Range("I31").Value = 33334552
$ws.Range("J31").Value = 25576.309
$ws.Range("K31").Value = 33334552
$ws.Range("L31").Value = 25576.309
$ws.Range("M31").Value = -33334257
$ws.Range("N31").Value = -26166.309

$ws.Range("H34").Value = 23264396
$ws.Range("I34").Value = 33334552
$ws.Range("J34").Value = 25576.309
$ws.Range("K34").Value = 33334552
$ws.Range("L34").Value = 25576.309
$ws.Range("M34").Value = -33334350
$ws.Range("N34").Value = -25980.309

$ws.Range("H132").Value = 7095950.5
$ws.Range("I132").Value = 9260944
$ws.Range("J132").Value = 10516.182
$ws.Range("K132").Value = 27782832
$ws.Range("L132").Value = 31548.546
$ws.Range("M132").Value = -27780302
$ws.Range("N132").Value = -36608.546

$ws.Range("H134").Value = 26786802
$ws.Range("I134").Value = 32052308
$ws.Range("J134").Value = 6251331.5
$ws.Range("K134").Value = 96156924
$ws.Range("L134").Value = 18753994.5
$ws.Range("M134").Value = -96154389
$ws.Range("N134").Value = -18759064.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 440.16666
$ws.Range("I5").Value = 371.0909
$ws.Range("J5").Value = 1200
$ws.Range("K5").Value = 1113.2727
$ws.Range("L5").Value = 3600
$ws.Range("M5").Value = -1001.2727
$ws.Range("N5").Value = -3824

$ws.Range("H122").Value = 718.75
$ws.Range("I122").Value = 706.931
$ws.Range("J122").Value = 833
$ws.Range("K122").Value = 6362.379000000001
$ws.Range("L122").Value = 7497
$ws.Range("M122").Value = -3912.379000000001
$ws.Range("N122").Value = -12397

$ws.Range("H135").Value = 440.16666
$ws.Range("I135").Value = 371.0909
$ws.Range("J135").Value = 1200
$ws.Range("K135").Value = 3339.8181
$ws.Range("L135").Value = 10800
$ws.Range("M135").Value = -804.8181
$ws.Range("N135").Value = -15870

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1092.3077
$ws.Range("I113").Value = 1157.375
$ws.Range("J113").Value = 988.2
$ws.Range("K113").Value = 1157.375
$ws.Range("L113").Value = 988.2
$ws.Range("M113").Value = 1012.625
$ws.Range("N113").Value = -5328.2

$ws.Range("H132").Value = 29414040
$ws.Range("I132").Value = 38462604
$ws.Range("J132").Value = 6210
$ws.Range("K132").Value = 115387812
$ws.Range("L132").Value = 18630
$ws.Range("M132").Value = -115385282
$ws.Range("N132").Value = -23690

$ws.Range("H141").Value = 46283.332
$ws.Range("J141").Value = 46283.332
$ws.Range("L141").Value = 46283.332
$ws.Range("N141").Value = -56643.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1597.3077
$ws.Range("I61").Value = 1085
$ws.Range("J61").Value = 2750
$ws.Range("K61").Value = 1085
$ws.Range("L61").Value = 2750
$ws.Range("M61").Value = -883
$ws.Range("N61").Value = -3154

$ws.Range("H113").Value = 1597.3077
$ws.Range("I113").Value = 1085
$ws.Range("J113").Value = 2750
$ws.Range("K113").Value = 1085
$ws.Range("L113").Value = 2750
$ws.Range("M113").Value = 1085
$ws.Range("N113").Value = -7090

$ws.Range("H132").Value = 2986250
$ws.Range("I132").Value = 4082547.8
$ws.Range("J132").Value = 1883.8889
$ws.Range("K132").Value = 12247643.4
$ws.Range("L132").Value = 5651.6667
$ws.Range("M132").Value = -12245113.4
$ws.Range("N132").Value = -10711.6667

$ws.Range("H136").Value = 5282.1665
$ws.Range("I136").Value = 5766.387
$ws.Range("J136").Value = 2280
$ws.Range("K136").Value = 17299.161
$ws.Range("L136").Value = 6840
$ws.Range("M136").Value = -14749.161
$ws.Range("N136").Value = -11940

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 20398
$ws.Range("J114").Value = 20398
$ws.Range("L114").Value = 20398
$ws.Range("N114").Value = -29076  # new cell added in target

$ws.Range("H132").Value = 2141130
$ws.Range("I132").Value = 3572199
$ws.Range("J132").Value = 319769.53
$ws.Range("K132").Value = 10716597
$ws.Range("L132").Value = 959308.5900000001
$ws.Range("M132").Value = -10714067
$ws.Range("N132").Value = -964368.5900000001

$ws.Range("H136").Value = 18633316
$ws.Range("I136").Value = 10203350
$ws.Range("K136").Value = 30610050
$ws.Range("M136").Value = -30607500

$ws.Range("H140").Value = 32601
$ws.Range("J140").Value = 32601
$ws.Range("L140").Value = 32601
$ws.Range("N140").Value = -42961

$ws.Range("H141").Value = 35000
$ws.Range("J141").Value = 35000
$ws.Range("L141").Value = 35000
$ws.Range("N141").Value = -45360
